$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2 through 302) holds a "Förändrad" (last-changed) date that
# was bumped by one day, from 2023-09-20 (serial 45189) to 2023-09-21
# (serial 45190), for every data row in the sheet.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 302) { $lastRow = 302 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45189) {
        $cell.Value2 = 45190
    }
}
